$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage (prevents Excel
# from auto-converting numeric-looking strings like '302.34' into numbers),
# and resets the cell style afterwards so formatting is left untouched.
function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "42.923.93"
$ws.Range("E2").Value = "  +0.14%  "
# Row 3
Set-TextCell "D3" "2.361.68"
$ws.Range("E3").Value = "  +1.60%  "
# Row 4
$ws.Range("E4").Value = "  +0.07%  "
# Row 5
Set-TextCell "D5" "302.34"
$ws.Range("E5").Value = "  +0.21%  "
# Row 6
Set-TextCell "D6" "95.43"
$ws.Range("E6").Value = "  -0.37%  "
# Row 7
$ws.Range("E7").Value = "  -0.04%  "
# Row 8
$ws.Range("E8").Value = "  -0.60%  "
# Row 9
Set-TextCell "D9" "0.482"
$ws.Range("E9").Value = "  -2.20%  "
# Row 10
Set-TextCell "D10" "34.01"
$ws.Range("E10").Value = "  -0.72%  "
# Row 11
$ws.Range("E11").Value = "  +3.36%  "
# Row 12
$ws.Range("E12").Value = "  -0.17%  "
# Row 13
Set-TextCell "D13" "18.37"
$ws.Range("E13").Value = "  -3.50%  "
# Row 14
Set-TextCell "D14" "2.729.63"
$ws.Range("E14").Value = "  +1.68%  "
# Row 15
Set-TextCell "D15" "6.68"
$ws.Range("E15").Value = "  -0.89%  "
# Row 16
Set-TextCell "D16" "2.374.79"
$ws.Range("E16").Value = "  +1.96%  "
# Row 17
$ws.Range("E17").Value = "  +0.02%  "
# Row 18
Set-TextCell "D18" "42.875.94"
$ws.Range("E18").Value = "  +0.15%  "
# Row 19
Set-TextCell "D19" "11.85"
$ws.Range("E19").Value = "  -2.85%  "
# Row 20
$ws.Range("E20").Value = "  +1.62%  "
# Row 21
Set-TextCell "D21" "0.0₃0884"
$ws.Range("E21").Value = "  -0.74%  "
# Row 22
Set-TextCell "D22" "67.96"
$ws.Range("E22").Value = "  +0.13%  "
# Row 23
Set-TextCell "D23" "234.95"
$ws.Range("E23").Value = "  -0.24%  "
# Row 24
$ws.Range("E24").Value = "  -4.01%  "
# Row 25
$ws.Range("E25").Value = "  -0.09%  "
# Row 26
$ws.Range("E26").Value = "  +0.49%  "
# Row 27
Set-TextCell "D27" "24.38"
$ws.Range("E27").Value = "  -0.20%  "
# Row 28
$ws.Range("E28").Value = "  +0.69%  "
# Row 29
$ws.Range("E29").Value = "  +1.77%  "
# Row 30
Set-TextCell "D30" "31.76"
$ws.Range("E30").Value = "  -1.47%  "
# Row 31
$ws.Range("E31").Value = "  +0.04%  "
# Row 32
$ws.Range("E32").Value = "  -0.23%  "
# Row 33
Set-TextCell "D33" "17.45"
$ws.Range("E33").Value = "  -2.32%  "
# Row 34
Set-TextCell "D34" "0.0721"
$ws.Range("E34").Value = "  +2.79%  "
# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D35" "1.84"
$ws.Range("E35").Value = "  +1.44%  "
# Row 36
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D36" "127.05"
$ws.Range("E36").Value = "  -12.83%  "
# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D37" "0.103"
$ws.Range("E37").Value = "  +3.40%  "
# Row 38
Set-TextCell "D38" "4.29"
$ws.Range("E38").Value = "  -2.40%  "
# Row 39
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D39" "2.28"
$ws.Range("E39").Value = "  -1.67%  "
# Row 40
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D40" "2.82"
$ws.Range("E40").Value = "  +2.40%  "
# Row 41
$ws.Range("E41").Value = "  -0.95%  "
# Row 42
Set-TextCell "D42" "21.23"
$ws.Range("E42").Value = "  -4.29%  "
# Row 43
Set-TextCell "D43" "1.927.52"
$ws.Range("E43").Value = "  +0.02%  "
# Row 44
$ws.Range("E44").Value = "  -0.35%  "
# Row 45
$ws.Range("E45").Value = "  +2.99%  "
# Row 46
$ws.Range("E46").Value = "  -9.42%  "
# Row 47
Set-TextCell "D47" "2.69"
$ws.Range("E47").Value = "  -1.98%  "
# Row 48
Set-TextCell "D48" "2.589.72"
$ws.Range("E48").Value = "  +1.44%  "
# Row 49
Set-TextCell "D49" "1.50"
$ws.Range("E49").Value = "  +1.53%  "
# Row 50
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell "D50" "71.49"
$ws.Range("E50").Value = "  -1.66%  "
# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D51" "1.14"
$ws.Range("E51").Value = "  +1.36%  "
